$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title in A1 to reflect the new migration wave date
$ws.Range("A1").Value = "Hotcarding Spreadsheet - Migration Wave 10/16/2025"

# Add the new data row (row 3)
$ws.Range("A3:J3").NumberFormat = "@"
$ws.Range("A3").Value = "10/16/2025"
$ws.Range("B3").Value = "YYY"
$ws.Range("C3").Value = "123ABX007"
$ws.Range("D3").Value = "FISB"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "PaymentsOne Debit"
$ws.Range("G3").Value = "Basic"
$ws.Range("H3").Value = "Offshore"
$ws.Range("I3").Value = "NA"
$ws.Range("J3").Value = "NA"
$ws.Range("A3:J3").Style = "Normal"
